$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# The "Name" row (row 4) loses its value; the "Title" row (row 5) now holds
# the text that used to sit in the Name row's value, and the old Title
# value is dropped entirely.
$wsMeta.Range("B4").Value = ""
$wsMeta.Range("B5").Value = "Mapping Métier/CDA/FHIR : ""Document Referencé"""

# Refresh the Date value
$wsMeta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# --- Mapping Table 0 sheet ---
$wsMap0 = $wb.Worksheets.Item("Mapping Table 0")
$wsMap0.Range("D8").Value = "FRCDADocumentAttache.component:frTypeDocumentAttache"

# --- Mapping Table 1 sheet ---
$wsMap1 = $wb.Worksheets.Item("Mapping Table 1")
$wsMap1.Range("A7").Value = "FRCDADocumentAttache.component:frTypeDocumentAttache"
